$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 363, shifting existing rows 363-433 down to 364-434
$ws.Rows.Item(363).Insert()

# Populate the newly inserted row 363 with the new data entry
$ws.Cells.Item(363, 1).Value = 5
$ws.Cells.Item(363, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(363, 3).Value = "Maule"
$ws.Cells.Item(363, 4).Value = 44637
$ws.Cells.Item(363, 5).Value = 7
$ws.Cells.Item(363, 6).Value = 100114001
$ws.Cells.Item(363, 7).Value = "Papa"
$ws.Cells.Item(363, 8).Value = "Rodeo"
$ws.Cells.Item(363, 9).Value = "1a (cosecha)"
$ws.Cells.Item(363, 10).Value = 1500
$ws.Cells.Item(363, 11).Value = 7000
$ws.Cells.Item(363, 12).Value = 7000
$ws.Cells.Item(363, 13).Value = 7000
$ws.Cells.Item(363, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(363, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(363, 16).Value = 280
$ws.Cells.Item(363, 17).Value = 25
$ws.Cells.Item(363, 18).Value = "Hortaliza"
